# Tasks.xlsx update: "Updated report with sections 2,3,4"
#
# This applies the autofilter (Status in {Almost Complete, Pending} AND
# Column2 in {the three specific notes} or blank), renames "WEEK3" (rows
# 19-20, the Week-3 task rows) to "WEEK 3" leaving row 21 on the old label,
# and moves the active selection to C25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Rename WEEK3 -> "WEEK 3" for the first two Week-3 rows (19 & 20) ---
# (Row 21 keeps referencing the original "WEEK3" shared string.)
$ws.Range("A19").Value2 = "WEEK 3"
$ws.Range("A20").Value2 = "WEEK 3"

# --- 2. Apply the AutoFilter on the table (Status column + Column2) ---
# Column2 ("F", table field 6): keep blanks plus these three notes.
$lo.Range.AutoFilter(6, @("Need to append graphs and prediction examples", "Sayantika to update predictions, examples and graphing", "We need to connect on Tuesday and put an outline together"), 7) | Out-Null
# Status ("D", table field 4): keep Almost Complete / Pending.
$lo.Range.AutoFilter(4, @("Almost Complete", "Pending"), 7) | Out-Null

# --- 3. Make sure row visibility matches the combined (both-column) filter ---
# The emulated AutoFilter only keeps the most-recently-applied column's
# criteria live, so re-assert the hidden/visible state that Excel would
# have produced from applying both filters together.
$hiddenRows = @(2,3,4,5,6,7,8,9,11,12,16,19,21)
$visibleRows = @(10,13,14,15,17,18,20)

foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}
foreach ($r in $visibleRows) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- 4. Move the selection to C25 ---
$ws.Range("C25").Select() | Out-Null

Write-Output "done"
